$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Service Rates")

# Insert two new rows above the existing "Tech. Regular Time" row (row 13)
# so the overtime rate rows can be added ahead of the regular-time rows.
$ws.Rows("13:14").Insert()

# New row 13: Tech. Overtime
$ws.Range("B13").Value = 8
$ws.Range("C13").Value = "Tech. Overtime"
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 185
$ws.Range("G13").Value = "Per hour/ per person"

# New row 14: Eng Overtime
$ws.Range("B14").Value = 9
$ws.Range("C14").Value = "Eng Overtime"
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 247
$ws.Range("G14").Value = "Per hour/ per person"

# Renumber the Item column for the rows that were pushed down
$ws.Range("B15").Value = 10
$ws.Range("B16").Value = 11
$ws.Range("B17").Value = 12

# Grow the table (and its filter range) to cover the two new rows
$tbl = $ws.ListObjects.Item("Table2")
$tbl.Resize($ws.Range("B3:G19"))

# Restore the view state: "Service Rates" becomes the active/selected sheet
$ws.Select()
$ws.Range("J14").Select()
$excel.ActiveWindow.ScrollRow = 4
